$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.793.55"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.755.97"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.16"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5051"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.20"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2658"
$ws.Range("E9").Value = "  +6.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06204"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.752.06"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.69"
$ws.Range("E12").Value = "  +6.31%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.06920"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5980"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.483"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.44"
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.820.83"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006814"
$ws.Range("E20").Value = "  +7.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.64"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.974.71"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.073"
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.254"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.57"
$ws.Range("E26").Value = "  +4.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.829"
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.450"
$ws.Range("E28").Value = "  +5.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.02"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.53"
$ws.Range("E30").Value = "  +3.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08186"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.663"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.383"
$ws.Range("E33").Value = "  +7.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04401"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9985"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.650"
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9990"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6004"
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.721"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01553"
$ws.Range("E40").Value = "  +6.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.936"
$ws.Range("E41").Value = "  -8.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.0000"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.05"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3795"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("E45").Value = "  -6.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.942"
$ws.Range("E46").Value = "  -4.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05492"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1092"
$ws.Range("E48").Value = "  +6.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.936"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.672"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "29.75"
$ws.Range("E51").Value = "  +1.96%  "
